$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 56, pushing existing row 56 (and below) down by one.
$ws.Rows("56:56").Insert()

# Populate the newly inserted row 56 with the new data record.
$ws.Range("A56").Value = 4
$ws.Range("B56").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C56").Value = "Los Lagos"
$ws.Range("D56").Value = 44985
$ws.Range("E56").Value = 10
$ws.Range("F56").Value = 100112043
$ws.Range("G56").Value = "Pepino dulce"
$ws.Range("H56").Value = "Cultivar IV Región"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 120
$ws.Range("K56").Value = 20000
$ws.Range("L56").Value = 22000
$ws.Range("M56").Value = 21000
$ws.Range("N56").Value = "$/bandeja 18 kilos"
$ws.Range("O56").Value = "Provincia de Limarí"
$ws.Range("P56").Value = 1167
$ws.Range("Q56").Value = 18
$ws.Range("R56").Value = "Hortaliza"
